$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -17.52953731255739
$ws.Range("C2").Value = 1.993559407765769
$ws.Range("D2").Value = -17.52953731255739
$ws.Range("E2").Value = -17.52953731255739
$ws.Range("F2").Value = -17.52953731255739
$ws.Range("G2").Value = -17.52953731255739
$ws.Range("H2").Value = -17.52953731255739
$ws.Range("I2").Value = -17.52953731255739
$ws.Range("J2").Value = -17.52953731255739
$ws.Range("K2").Value = -17.52953731255739
$ws.Range("B3").Value = -17.52953731255739
$ws.Range("C3").Value = -17.52953731255739
$ws.Range("D3").Value = -17.52953731255739
$ws.Range("E3").Value = -17.52953731255739
$ws.Range("F3").Value = -17.52953731255739
$ws.Range("G3").Value = -17.52953731255739
$ws.Range("H3").Value = -17.52953731255739
$ws.Range("I3").Value = 2.998280419409443
$ws.Range("J3").Value = -17.52953731255739
$ws.Range("K3").Value = -17.52953731255739
$ws.Range("B4").Value = -17.52953731255739
$ws.Range("C4").Value = 2.184173066206692
$ws.Range("D4").Value = 2.258341340269291
$ws.Range("E4").Value = -17.52953731255739
$ws.Range("F4").Value = 3.366216195313942
$ws.Range("G4").Value = -17.52953731255739
$ws.Range("H4").Value = 1.487241056477941
$ws.Range("I4").Value = -17.52953731255739
$ws.Range("J4").Value = -17.52953731255739
$ws.Range("K4").Value = -17.52953731255739
$ws.Range("B5").Value = -17.52953731255739
$ws.Range("C5").Value = 1.87129418759832
$ws.Range("D5").Value = -17.52953731255739
$ws.Range("E5").Value = -17.52953731255739
$ws.Range("F5").Value = -17.52953731255739
$ws.Range("G5").Value = 3.119751323268142
$ws.Range("H5").Value = -17.52953731255739
$ws.Range("I5").Value = -17.52953731255739
$ws.Range("J5").Value = -17.52953731255739
$ws.Range("K5").Value = -17.52953731255739
$ws.Range("B6").Value = -17.52953731255739
$ws.Range("C6").Value = -17.52953731255739
$ws.Range("D6").Value = -17.52953731255739
$ws.Range("E6").Value = -17.52953731255739
$ws.Range("F6").Value = -17.52953731255739
$ws.Range("G6").Value = -17.52953731255739
$ws.Range("H6").Value = -17.52953731255739
$ws.Range("I6").Value = -17.52953731255739
$ws.Range("J6").Value = -17.52953731255739
$ws.Range("K6").Value = -17.52953731255739
$ws.Range("B7").Value = 2.65665069487889
$ws.Range("C7").Value = -17.52953731255739
$ws.Range("D7").Value = -17.52953731255739
$ws.Range("E7").Value = -17.52953731255739
$ws.Range("F7").Value = -17.52953731255739
$ws.Range("G7").Value = -17.52953731255739
$ws.Range("H7").Value = -17.52953731255739
$ws.Range("I7").Value = -17.52953731255739
$ws.Range("J7").Value = -17.52953731255739
$ws.Range("K7").Value = -17.52953731255739
$ws.Range("B8").Value = -17.52953731255739
$ws.Range("C8").Value = -17.52953731255739
$ws.Range("D8").Value = -17.52953731255739
$ws.Range("E8").Value = 1.794768511409633
$ws.Range("F8").Value = -17.52953731255739
$ws.Range("G8").Value = -17.52953731255739
$ws.Range("H8").Value = -17.52953731255739
$ws.Range("I8").Value = -17.52953731255739
$ws.Range("J8").Value = -17.52953731255739
$ws.Range("K8").Value = -17.52953731255739
$ws.Range("B9").Value = 3.775496283614083
$ws.Range("C9").Value = -17.52953731255739
$ws.Range("D9").Value = -17.52953731255739
$ws.Range("E9").Value = -17.52953731255739
$ws.Range("F9").Value = -17.52953731255739
$ws.Range("G9").Value = -17.52953731255739
$ws.Range("H9").Value = -17.52953731255739
$ws.Range("I9").Value = -17.52953731255739
$ws.Range("J9").Value = -17.52953731255739
$ws.Range("K9").Value = -17.52953731255739
$ws.Range("B10").Value = -17.52953731255739
$ws.Range("C10").Value = -17.52953731255739
$ws.Range("D10").Value = -17.52953731255739
$ws.Range("E10").Value = -17.52953731255739
$ws.Range("F10").Value = -17.52953731255739
$ws.Range("G10").Value = -17.52953731255739
$ws.Range("H10").Value = -17.52953731255739
$ws.Range("I10").Value = 1.353084842218539
$ws.Range("J10").Value = -17.52953731255739
$ws.Range("K10").Value = 1.545871325827293
$ws.Range("B11").Value = -17.52953731255739
$ws.Range("C11").Value = -17.52953731255739
$ws.Range("D11").Value = -17.52953731255739
$ws.Range("E11").Value = 2.861547860754214
$ws.Range("F11").Value = -17.52953731255739
$ws.Range("G11").Value = 2.386629356232455
$ws.Range("H11").Value = -17.52953731255739
$ws.Range("I11").Value = -17.52953731255739
$ws.Range("J11").Value = -17.52953731255739
$ws.Range("K11").Value = 1.533702388344111
$ws.Range("B12").Value = -17.52953731255739
$ws.Range("C12").Value = -17.52953731255739
$ws.Range("D12").Value = -17.52953731255739
$ws.Range("E12").Value = -17.52953731255739
$ws.Range("F12").Value = -17.52953731255739
$ws.Range("G12").Value = -17.52953731255739
$ws.Range("H12").Value = -17.52953731255739
$ws.Range("I12").Value = -17.52953731255739
$ws.Range("J12").Value = -17.52953731255739
$ws.Range("K12").Value = -17.52953731255739
$ws.Range("B13").Value = -17.52953731255739
$ws.Range("C13").Value = -17.52953731255739
$ws.Range("D13").Value = -17.52953731255739
$ws.Range("E13").Value = 2.393437546561615
$ws.Range("F13").Value = -17.52953731255739
$ws.Range("G13").Value = -17.52953731255739
$ws.Range("H13").Value = -17.52953731255739
$ws.Range("I13").Value = -17.52953731255739
$ws.Range("J13").Value = -17.52953731255739
$ws.Range("K13").Value = 1.825557796359109
$ws.Range("B14").Value = -17.52953731255739
$ws.Range("C14").Value = -17.52953731255739
$ws.Range("D14").Value = 1.243991498023847
$ws.Range("E14").Value = -17.52953731255739
$ws.Range("F14").Value = -17.52953731255739
$ws.Range("G14").Value = -17.52953731255739
$ws.Range("H14").Value = -17.52953731255739
$ws.Range("I14").Value = -17.52953731255739
$ws.Range("J14").Value = -17.52953731255739
$ws.Range("K14").Value = 2.257161457628485
$ws.Range("B15").Value = -17.52953731255739
$ws.Range("C15").Value = -17.52953731255739
$ws.Range("D15").Value = 1.257705280184706
$ws.Range("E15").Value = -17.52953731255739
$ws.Range("F15").Value = -17.52953731255739
$ws.Range("G15").Value = -17.52953731255739
$ws.Range("H15").Value = -17.52953731255739
$ws.Range("I15").Value = -17.52953731255739
$ws.Range("J15").Value = -17.52953731255739
$ws.Range("K15").Value = -17.52953731255739
$ws.Range("B16").Value = -17.52953731255739
$ws.Range("C16").Value = -17.52953731255739
$ws.Range("D16").Value = -17.52953731255739
$ws.Range("E16").Value = -17.52953731255739
$ws.Range("F16").Value = -17.52953731255739
$ws.Range("G16").Value = -17.52953731255739
$ws.Range("H16").Value = -17.52953731255739
$ws.Range("I16").Value = -17.52953731255739
$ws.Range("J16").Value = 4.321920850816444
$ws.Range("K16").Value = -17.52953731255739
$ws.Range("B17").Value = -17.52953731255739
$ws.Range("C17").Value = 1.882877476552768
$ws.Range("D17").Value = 2.158543920843704
$ws.Range("E17").Value = -17.52953731255739
$ws.Range("F17").Value = -17.52953731255739
$ws.Range("G17").Value = -17.52953731255739
$ws.Range("H17").Value = 0.6865501279179931
$ws.Range("I17").Value = 1.311845826996312
$ws.Range("J17").Value = -17.52953731255739
$ws.Range("K17").Value = -17.52953731255739
$ws.Range("B18").Value = -17.52953731255739
$ws.Range("C18").Value = -17.52953731255739
$ws.Range("D18").Value = -17.52953731255739
$ws.Range("E18").Value = -17.52953731255739
$ws.Range("F18").Value = -17.52953731255739
$ws.Range("G18").Value = -17.52953731255739
$ws.Range("H18").Value = 0.9688016074005524
$ws.Range("I18").Value = 0.8208192121333425
$ws.Range("J18").Value = -17.52953731255739
$ws.Range("K18").Value = -17.52953731255739
$ws.Range("B19").Value = -17.52953731255739
$ws.Range("C19").Value = -17.52953731255739
$ws.Range("D19").Value = 1.634301512182141
$ws.Range("E19").Value = -17.52953731255739
$ws.Range("F19").Value = -17.52953731255739
$ws.Range("G19").Value = -17.52953731255739
$ws.Range("H19").Value = 1.739619050993904
$ws.Range("I19").Value = 1.768464825284994
$ws.Range("J19").Value = -17.52953731255739
$ws.Range("K19").Value = -17.52953731255739
$ws.Range("B20").Value = -17.52953731255739
$ws.Range("C20").Value = 0.9532903349066769
$ws.Range("D20").Value = 1.529482289819035
$ws.Range("E20").Value = -17.52953731255739
$ws.Range("F20").Value = 3.276223088281324
$ws.Range("G20").Value = -17.52953731255739
$ws.Range("H20").Value = 2.274160262629553
$ws.Range("I20").Value = 0.8471688982863174
$ws.Range("J20").Value = -17.52953731255739
$ws.Range("K20").Value = 2.550875576350224
$ws.Range("B21").Value = -17.52953731255739
$ws.Range("C21").Value = 1.130991383643336
$ws.Range("D21").Value = -17.52953731255739
$ws.Range("E21").Value = 2.002987934985138
$ws.Range("F21").Value = -17.52953731255739
$ws.Range("G21").Value = 2.603641281908413
$ws.Range("H21").Value = 2.447084143869486
$ws.Range("I21").Value = -17.52953731255739
$ws.Range("J21").Value = -17.52953731255739
$ws.Range("K21").Value = -17.52953731255739
